# Add "NA" under duplicate_image_filename (column E) for data rows 2-21
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2:E21").Value = "NA"
